$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New data rows to append after existing table data (rows 16-19: LRelu agents; rows 20-23: Elu agents, partial)
# Note: scientific-notation literals (e.g. 1.0E-2) are not parsed by this PowerShell engine,
# so equivalent plain-decimal literals are used instead (they round-trip to the same double).
$data = @(
    @("Raw-Sutton-40-lrelu",    "Lrelu", "Sutton",    319, 0.010033394684030699, 352, 0.0311701877751514,  308, 0.0178892986171405),
    @("Raw-Tesauro89-40-lrelu", "Lrelu", "Tesauro89", 165, 0.0085945415559447805, 187, 0.0144900887994673,  165, 0.011670147232056199),
    @("Raw-Tesauro92-40-lrelu", "Lrelu", "Tesauro92", 165, 0.0091649879250330201, 254, 0.014580407029913801, 198, 0.0138817468289366),
    @("Raw-GnuBg-40-lrelu",     "Lrelu", "GnuBg",     220, 0.00919357906345529,   253, 0.014600829719728,   143, 0.017127940607985)
)

$startRow = 16
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 9).Value = $row[8]
}

# Elu rows: only Agent Name on the first (Sutton) row, and Activation/Codec on all four
$eluRows = @(
    @("Raw-Sutton-40-elu", "Elu", "Sutton"),
    @($null,               "Elu", "Tesauro89"),
    @($null,               "Elu", "Tesauro92"),
    @($null,               "Elu", "GnuBg")
)

$eluStart = 20
for ($i = 0; $i -lt $eluRows.Count; $i++) {
    $r = $eluStart + $i
    $row = $eluRows[$i]
    if ($row[0] -ne $null) {
        $ws.Cells.Item($r, 1).Value = $row[0]
    }
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}

# Resize the table to include the new rows (expands ref and autoFilter to A3:I23)
$table = $ws.ListObjects.Item("Table2")
$table.Resize($ws.Range("A3:I23"))


# Update the selection to match the final cursor position
$ws.Range("A21").Select()
